# Atualizações na coluna existentes planilha BMA
# - Add a new worksheet "Planilha2" with a small helper calculation
#   (B1 = 4, A4:A6 random numbers bounded by B1, C4:C7 running "SMALL"
#   differences that always sum back up to B1).
# - Make the new sheet the active sheet/tab.
# - Update the selection on the original "Necessidades" sheet to F2.

$wb = $excel.ActiveWorkbook

# --- Work on the existing "Necessidades" sheet first -----------------
$wsNecessidades = $wb.Worksheets.Item("Necessidades")
[void]$wsNecessidades.Select()
[void]$wsNecessidades.Range("F2").Select()

# --- Add the new sheet -------------------------------------------------
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNecessidades)
$wsNew.Name = "Planilha2"

# Seed value
$wsNew.Range("B1").Value = 4

# Random values bounded by B1 (volatile formulas)
$wsNew.Range("A4").Formula = "=RANDBETWEEN(0, `$B`$1)"
$wsNew.Range("A5").Formula = "=RANDBETWEEN(0, `$B`$1)"
$wsNew.Range("A6").Formula = "=RANDBETWEEN(0, `$B`$1)"

# Running "SMALL" based split of B1 across four buckets
$wsNew.Range("C4").Formula = "=SMALL(`$A`$4:`$A`$6, 1)"
$wsNew.Range("C5").Formula = "=SMALL(`$A`$4:`$A`$6, 2) - SMALL(`$A`$4:`$A`$6, 1)"
$wsNew.Range("C6").Formula = "=SMALL(`$A`$4:`$A`$6, 3) - SMALL(`$A`$4:`$A`$6, 2)"
$wsNew.Range("C7").Formula = "=`$B`$1 - SMALL(`$A`$4:`$A`$6, 3)"

[void]$wb.Application.Calculate()

# Selection on the new sheet
[void]$wsNew.Range("C10").Select()

# Make sure the new sheet ends up as the active/selected tab
[void]$wsNew.Select()
